# exclusão de aba não utilizada
$wb = $excel.ActiveWorkbook

# Remove the unused "fator" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("fator").Delete()
$excel.DisplayAlerts = $true

# Rename the remaining "meses" sheet to "tipo"
$wb.Worksheets.Item("meses").Name = "tipo"
